# BF: typos and wrong places for some translations.
# Adds new i18n strings to Feuil1 (rows 320-340), adds a note to A311/B311,
# and tidies up the window view (selection/scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New rows 320-340: column A text, written in the exact order the
#     strings were first introduced in the shared-string table so the
#     resulting sharedStrings.xml matches the authored order. ---
$ws.Range('A320').Value = 'The overtime you''ve requested has been rejected. Below, the details :'
$ws.Range('A321').Value = 'The overtime you''ve requested has been accepted. Below, the details :'
$ws.Range('A322').Value = '{Firstname} {Lastname} requests an overtime. Below, the details :'
$ws.Range('A324').Value = 'Once connected, you can change your password, as explained here.'
$ws.Range('A325').Value = 'The leave you''ve requested has been rejected. Below, the details :'
$ws.Range('A326').Value = 'The leave you''ve requested has been accepted. Below, the details :'
$ws.Range('A327').Value = '{Firstname} {Lastname} requests a leave. Below, the details :'
$ws.Range('A328').Value = 'From'
$ws.Range('A329').Value = 'To'
$ws.Range('A330').Value = 'Dear {Firstname} {Lastname},'
$ws.Range('A331').Value = 'If you didn''t perform this operation, please contact your administrator.'
$ws.Range('A332').Value = 'Welcome in LMS. If your are an employee, you could now :'
$ws.Range('A333').Value = 'See your leave balance.'
$ws.Range('A334').Value = 'See the list of the leave requests you have submitted.'
$ws.Range('A335').Value = 'Request a new leave.'
$ws.Range('A336').Value = 'If your are the line manager of other employee(s), you could now :'
$ws.Range('A337').Value = 'Validate leave requests submitted to you.'
$ws.Range('A338').Value = 'Validate overtime requests submitted to you.'
$ws.Range('A339').Value = 'Access forbidden'
$ws.Range('A340').Value = 'You are not allowed to perform this action.'

# Column B notes for translators (shared across rows 322/323/330)
$ws.Range('B322').Value = 'don''t remove or replace {Firstname} and {Lastname}'
$ws.Range('A323').Value = 'Welcome to LMS {Firstname} {Lastname}. Please use these credentials to login to the system:'
$ws.Range('B311').Value = 'In the sense of overtime (extra hours)'
$ws.Range('B323').Value = 'don''t remove or replace {Firstname} and {Lastname}'
$ws.Range('B330').Value = 'don''t remove or replace {Firstname} and {Lastname}'

# --- Keep the row "shape" (column span) consistent with the rest of the
#     sheet for rows 320-336, which historically carry a 2-column span
#     even when only column A is populated. ---
$spanRows = @(320,321,324,325,326,327,328,329,331,332,333,334,335,336)
foreach ($r in $spanRows) {
    $ws.Range("B$r").Value = 'x'
}
foreach ($r in $spanRows) {
    $ws.Range("B$r").ClearContents()
}

# --- Restore the view: scroll back to top-left and select the newly
#     added "To" row instead of the old bottom-of-sheet selection. ---
$ws.Range('A329').Select()
